$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 42-43; existing rows 42-57 shift down to 44-59.
$ws.Rows("42:43").Insert()

# New row 42: Jengibre, Primera, week of 2021-09-27 (serial 44466)
$ws.Cells.Item(42, 1).Value = 9
$ws.Cells.Item(42, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(42, 3).Value = "Metropolitana"
$ws.Cells.Item(42, 4).Value = 44466
$ws.Cells.Item(42, 5).Value = 13
$ws.Cells.Item(42, 6).Value = 100114007
$ws.Cells.Item(42, 7).Value = "Jengibre"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 1150
$ws.Cells.Item(42, 11).Value = 14000
$ws.Cells.Item(42, 12).Value = 15000
$ws.Cells.Item(42, 13).Value = 14500
$ws.Cells.Item(42, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(42, 15).Value = "Perú"
$ws.Cells.Item(42, 16).Value = 1115
$ws.Cells.Item(42, 17).Value = 13
$ws.Cells.Item(42, 18).Value = "Hortaliza"

# New row 43: Jengibre, Segunda, week of 2021-09-27 (serial 44466)
$ws.Cells.Item(43, 1).Value = 9
$ws.Cells.Item(43, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(43, 3).Value = "Metropolitana"
$ws.Cells.Item(43, 4).Value = 44466
$ws.Cells.Item(43, 5).Value = 13
$ws.Cells.Item(43, 6).Value = 100114007
$ws.Cells.Item(43, 7).Value = "Jengibre"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Segunda"
$ws.Cells.Item(43, 10).Value = 790
$ws.Cells.Item(43, 11).Value = 12000
$ws.Cells.Item(43, 12).Value = 12000
$ws.Cells.Item(43, 13).Value = 12000
$ws.Cells.Item(43, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(43, 15).Value = "Perú"
$ws.Cells.Item(43, 16).Value = 923
$ws.Cells.Item(43, 17).Value = 13
$ws.Cells.Item(43, 18).Value = "Hortaliza"

# Apply the date number format (matching the other date cells in column D)
$ws.Range("D42:D43").NumberFormat = $ws.Range("D44").NumberFormat
